# Commit: "changing first sheet and adding code"
# Adds a new worksheet "Resumo_Cores" right after the existing "Item" sheet
# and populates its header row with the color/wire-gauge summary columns.

$wb = $excel.ActiveWorkbook

# The existing (first) sheet - content/formatting stays as-is.
$itemSheet = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after "Item" so sheet order is preserved.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $itemSheet)
$newSheet.Name = "Resumo_Cores"

$newSheet.Range("A1").Value = "Cor"
$newSheet.Range("B1").Value = "UL STYLE"
$newSheet.Range("C1").Value = "WIRE GAUGE"
$newSheet.Range("D1").Value = "TEMP RATING"

# Keep "Item" as the selected/active tab, matching the original workbook.
$itemSheet.Activate()
